$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2's target cluster changes from "ECs" to "MuSCs" (the old "ECs" shared
# string is dropped and the values previously on row 3 are folded into row 2)
$ws.Range("D2").Value = "MuSCs"

# Updated aggregate values for row 2 (K2:T2)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.4426103333333333
$ws.Range("N2").Value = 1.327831
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.07152582986666667
$ws.Range("R2").Value = 0.6437324688
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1

# Row 3 (the old separate MuSCs row) is no longer needed; remove it entirely
$ws.Rows(3).Delete()
